# This script re-orders the content of several data rows in the "Artfynd"
# sheet. The workbook records individual species observations, one per row
# (columns A:AY). A set of rows need their entire row-content rotated among
# each other (the row *positions* stay the same; the *data* that lives in
# them moves), per these cycles (content moves from each row to the next,
# wrapping at the end of the list):
#   3 -> 4 -> 7 -> 6 -> 3
#   10 -> 11 -> 12 -> 13 -> 10
#   23 -> 24 -> 23
#   37 -> 39 -> 37
#   38 -> 40 -> 38

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full data extent is A1:AY41 -> data columns are A (1) through AY (51).
$firstCol = 1
$lastCol = 51

function Get-RowValues($row) {
    $rng = $ws.Range($ws.Cells.Item($row, $firstCol), $ws.Cells.Item($row, $lastCol))
    return $rng.Value()
}

function Set-RowValues($row, $values) {
    $rng = $ws.Range($ws.Cells.Item($row, $firstCol), $ws.Cells.Item($row, $lastCol))
    $rng.Value = $values
}

$cycles = @(
    @(3, 4, 7, 6),
    @(10, 11, 12, 13),
    @(23, 24),
    @(37, 39),
    @(38, 40)
)

foreach ($cycle in $cycles) {
    # Snapshot every row in the cycle first so writes never clobber a row
    # before it has been read.
    $snapshots = @{}
    foreach ($row in $cycle) {
        $snapshots[$row] = Get-RowValues $row
    }

    $count = $cycle.Count
    for ($i = 0; $i -lt $count; $i++) {
        $srcRow = $cycle[$i]
        $dstRow = $cycle[($i + 1) % $count]
        Set-RowValues $dstRow $snapshots[$srcRow]
    }
}
